$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (formulas will recalc automatically)
$ws.Range("D5").Value = 1600
$ws.Range("E8").Value = 1850

# Update selection to match the new active cell/selection in sheetView
$ws.Range("F12").Select()
